$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Notified Production (MW)" values for rows 2..97 (column B)
$newB = @(249.749,246.398,242.562,239.851,232.206,225.794,221.343,215.796,203.927,200.271,195.829,190.87,184.661,179.079,174.598,170.409,165.129,162.833,159.561,156.313,146.621,144.578,142.626,140.528,138.326,138.975,139.92,141.167,130.7,132.828,134.899,136.978,132.843,131.746,131.627,131.184,140.881,139.694,139.847,139.836,159.441,159.49,159.711,158.988,161.536,162.513,162.288,163.184,167.535,171.695,176.141,180.376,191.067,202.355,212.429,223.967,244.759,261.394,276.635,293.032,315.516,327.263,338.597,350.391,383.797,404.337,424.162,444.898,493.061,513.91,536.133,558.796,593.374,603.645,613.691,625.3920000000001,641.533,647.958,654.026,660.431,668.573,672.7329999999999,675.629,678.681,680.313,684.021,688.458,691.741,691.1130000000001,693.045,694.575,696.1900000000001,0,0,0,0)

# Shift the timestamp column (A) forward by 2 days for rows 2..97,
# and write the new values into column B.
for ($i = 0; $i -lt $newB.Length; $i++) {
    $row = $i + 2
    $aCell = $ws.Cells.Item($row, 1)
    $aCell.Value2 = $aCell.Value2 + 2
    $ws.Cells.Item($row, 2).Value2 = $newB[$i]
}
